# feat: add 2022-Q1 data
#
# The workbook currently has two sheets: "2021-Q4" (fund holding detail)
# and "总计" (quarter-over-quarter summary).
# This script:
#   1. Inserts a new worksheet "2022-Q1" between "2021-Q4" and "总计",
#      populated with the new quarter's fund holding detail.
#   2. Inserts a new leading row into "总计" summarizing the 2022-Q1 data,
#      pushing the existing 2021-Q4 summary row down.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)      # "2021-Q4" - used as a style template
$total = $wb.Worksheets.Item(2)    # "总计"

# 1. Insert new worksheet "2022-Q1" right after "2021-Q4" (i.e. before "总计")
$newSheet = $wb.Worksheets.Add($total)
$newSheet.Name = "2022-Q1"

# NOTE: $total was obtained via a positional index lookup, which in this
# COM layer re-resolves live by index rather than holding a fixed
# reference. Since the new sheet was inserted before it, index 2 now
# refers to the new "2022-Q1" sheet instead of "总计". Re-fetch "总计"
# by name so subsequent edits target the right sheet.
$total = $wb.Worksheets.Item("总计")

# 2. Reuse the exact same cell styles as "2021-Q4" (bold + bordered +
#    centered header row, and bold + bordered + centered index column)
#    by copying formats from it instead of re-declaring new styles.
$ws1.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$ws1.Range("A2:A5").Copy()
$newSheet.Range("A2:A5").PasteSpecial(-4122)

# 3. Header row text
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# 4. Index column A
$newSheet.Range("A2").Value = 0
$newSheet.Range("A3").Value = 1
$newSheet.Range("A4").Value = 2
$newSheet.Range("A5").Value = 3

# 5. Row 2 - fund 506007
#    (leading "'" forces text storage so values like "506007"/"5.33"
#    keep their exact printed form instead of becoming numbers)
$newSheet.Range("B2").Value = "'506007"
$newSheet.Range("C2").Value = "广发科创板两年定开混合"
$newSheet.Range("D2").Value = "'5.33"
$newSheet.Range("E2").Value = "'92.08"
$newSheet.Range("F2").Value = "'4.73"
$newSheet.Range("G2").Value = "'0.2521"
$newSheet.Range("H2").Value = 5

# 6. Row 3 - fund 519951
$newSheet.Range("B3").Value = "'519951"
$newSheet.Range("C3").Value = "长信利泰灵活配置混合A"
$newSheet.Range("D3").Value = "'0.07"
$newSheet.Range("E3").Value = "'25.18"
$newSheet.Range("F3").Value = "'3.46"
$newSheet.Range("G3").Value = "'0.0024"
$newSheet.Range("H3").Value = 1

# 7. Row 4 - fund 008071
$newSheet.Range("B4").Value = "'008071"
$newSheet.Range("C4").Value = "长信利泰灵活配置混合E"
$newSheet.Range("D4").Value = "'0.01"
$newSheet.Range("E4").Value = "'25.18"
$newSheet.Range("F4").Value = "'3.46"
$newSheet.Range("G4").Value = "'0.0003"
$newSheet.Range("H4").Value = 1

# 8. Row 5 - fund 007863 (G5 is genuinely numeric 0, unlike the other
#    "0.xxxx" values in column G which are kept as text)
$newSheet.Range("B5").Value = "'007863"
$newSheet.Range("C5").Value = "长信利泰灵活配置混合C"
$newSheet.Range("D5").Value = "'0.00"
$newSheet.Range("E5").Value = "'25.18"
$newSheet.Range("F5").Value = "'3.46"
$newSheet.Range("G5").Value = 0
$newSheet.Range("H5").Value = 1

# 9. Update "总计" sheet: insert a new row 2 holding the 2022-Q1 summary,
#    pushing the existing 2021-Q4 summary row down to row 3.
$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

# Reuse the index-column style (bold + bordered + centered) from row 3
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.25

$total.Range("A3").Value = 1
